# Apply "All Unit Tests Added" commit changes (reverse-applied here: remove
# the CARD_006 test row that was present, update the "Last Runtime" dates,
# narrow column B, and move the active selection) to TestCases.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the CARD_006 test row contents and formatting (A7:E7) - row 7
# becomes an empty row like the rest of the sheet (F7 keeps its style-only cell)
$ws.Range("A7:E7").Clear()

# Shrink Table1 so it no longer includes the removed CARD_006 row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E6"))

# Update "Last Runtime" dates for remaining rows 2-6 from 44414 to 44413
$ws.Range("E2:E6").Value = 44413

# Narrow column B width (target stored width 52.8984375 pts ~= 52.14 chars)
$ws.Columns.Item(2).ColumnWidth = 52.14

# Move active cell selection to C6
$ws.Range("C6").Select()
